# "July Sapflow Plot level"
# Adds a new "Plot Level Sapflow" block to the bottom of the Main sheet:
#   - row 29: blank spacer row (same formatting as the existing spacer at row 25)
#   - row 30: A = "Plot Level Sapflow", B = "Calculated sectional area sums across plots"
#   - row 31: B = "July Sectional Sapflow (*sapflux/10000)"
#   - row 32: B = "Graphs"
# and updates the sheet's scroll/selection position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Recreate the blank "divider" row (row 29) using the same look as the
# existing divider row (row 25: vertical-center, no-wrap, customFormat).
$ws.Rows.Item(25).Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null

# New "Plot Level Sapflow" block.
$ws.Cells.Item(30, 1).Value = "Plot Level Sapflow"
$ws.Cells.Item(30, 2).Value = "Calculated sectional area sums across plots"
$ws.Cells.Item(31, 2).Value = "July Sectional Sapflow (*sapflux/10000)"
$ws.Cells.Item(32, 2).Value = "Graphs"

# Scroll the view down and move the selection to where the new rows were
# added (matches the saved view state: topLeftCell A20 / selection D35).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D35").Select() | Out-Null
